$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date serial value from 45208 to 45212 for rows 2-8
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45212
}
